$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original data had two rows: row 2 = FAPs/Rspo4/Lgr6 -> ECs, row 3 =
# FAPs/Rspo4/Lgr6 -> MuSCs. The ECs row is dropped entirely (with its now-unused
# "ECs" label pruned from the shared-string table), and the MuSCs row moves up
# to become the sole data row, with its specificity columns (O,P,Q,R,S,T)
# refreshed to newly recomputed TPM-based values.

$ws.Rows("2:2").Delete()

$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 0.047982992647
$ws.Range("R2").Value2 = 0.431846933823
$ws.Range("S2").Value2 = 1
$ws.Range("T2").Value2 = 1
